$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point (2026/01/27, 火, 8, 201) was logged. It belongs right
# after the existing 2026/01/27 rows (723/724 -> times 1,5), so insert a
# fresh row at 723 which pushes the rest of the table down by one, then
# populate it.
$ws.Rows.Item(723).Insert()

$ws.Cells.Item(723, 1).Value = "'2026/01/27"
$ws.Cells.Item(723, 1).Style = "Normal"
$ws.Cells.Item(723, 2).Value = "火"
$ws.Cells.Item(723, 3).Value = 8
$ws.Cells.Item(723, 4).Value = 201
